# CU-21 Consultar monedero - corrige la descripcion y el diagrama de robustez
#
# 1) "El sistema cierra la ventana "ConsultarMonederoView"." ->
#    "El sistema cierra la ventana "ConsultarMonederoView" y muestra la
#     ventana "VerMonederosView"."
# 2) "El sistema cierra las ventanas "ErrorView" y "ConsultarMonederoView"." ->
#    "El sistema cierra las ventanas "ErrorView", "ConsultarMonederoView" y
#     "VerMonederosView"."

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# The literal text "onsultarMonederoView" also occurs (split over runs as
# "Consultar"+"Monedero"+"View") earlier in the normal flow, so scope the
# search to start after the unique phrase "cierra la ventana " to land on
# the right sentence, then replace just the single first match.
$anchor1 = $d.Content.Duplicate
$anchor1.Find.Execute("cierra la ventana “C")
$scope1 = $d.Range($anchor1.Start, $d.Content.End)
$scope1.Find.Execute(
    "onsultarMonederoView”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "onsultarMonederoView” y muestra la ventana “VerMonederosView”.",
    1)

# --- Edit 2 -----------------------------------------------------------
$d.Content.Find.Execute(
    " y “ConsultarMonederoView”.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", “ConsultarMonederoView” y “VerMonederosView”.",
    1)
